$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '42.040.65'
$ws.Range("E2").Value = '  -0.19%  '

$ws.Range("D3").Value = '2.223.30'
$ws.Range("E3").Value = '  -0.88%  '

$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").Value = '243.58'
$ws.Range("E5").Value = '  -1.43%  '

$ws.Range("D6").Value = '0.628'
$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("D7").Value = '73.70'
$ws.Range("E7").Value = '  -0.61%  '

$ws.Range("E8").Value = '  +0.19%  '

$ws.Range("E9").Value = '  -0.33%  '

$ws.Range("D10").Value = '43.31'
$ws.Range("E10").Value = '  +5.90%  '

$ws.Range("D11").Value = '0.0961'
$ws.Range("E11").Value = '  +2.99%  '

$ws.Range("D12").Value = '7.11'
$ws.Range("E12").Value = '  +0.36%  '

$ws.Range("E13").Value = '  +0.46%  '

$ws.Range("E14").Value = '  -1.04%  '

$ws.Range("D15").Value = '0.845'
$ws.Range("E15").Value = '  -0.89%  '

$ws.Range("D16").Value = '2.238.69'
$ws.Range("E16").Value = '  +0.06%  '

$ws.Range("D17").Value = '41.955.69'
$ws.Range("E17").Value = '  -0.11%  '

$ws.Range("D18").Value = '0.0000109'
$ws.Range("E18").Value = '  +11.99%  '

$ws.Range("D19").Value = '6.21'
$ws.Range("E19").Value = '  +1.66%  '

$ws.Range("D20").Value = '72.23'
$ws.Range("E20").Value = '  +0.63%  '

$ws.Range("D21").Value = '10.17'
$ws.Range("E21").Value = '  +32.46%  '

$ws.Range("D22").Value = '229.56'
$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("E23").Value = '  -8.10%  '

$ws.Range("D24").Value = '11.60'
$ws.Range("E24").Value = '  +5.27%  '

$ws.Range("E25").Value = '  +0.11%  '

$ws.Range("E26").Value = '  +1.07%  '

$ws.Range("D27").Value = '2.28'
$ws.Range("E27").Value = '  -0.74%  '

$ws.Range("E28").Value = '  -1.67%  '

$ws.Range("D29").Value = '166.88'
$ws.Range("E29").Value = '  -2.77%  '

$ws.Range("D30").Value = '20.64'
$ws.Range("E30").Value = '  -0.03%  '

$ws.Range("D31").Value = '5.60'
$ws.Range("E31").Value = '  +13.97%  '

$ws.Range("D32").Value = '0.0798'
$ws.Range("E32").Value = '  -3.33%  '

$ws.Range("E33").Value = '  +0.90%  '

$ws.Range("D34").Value = '29.48'
$ws.Range("E34").Value = '  -1.86%  '

$ws.Range("E35").Value = '  -3.37%  '

$ws.Range("D36").Value = '4.31'
$ws.Range("E36").Value = '  -4.39%  '

$ws.Range("E37").Value = '  +0.57%  '

$ws.Range("D38").Value = '13.04'
$ws.Range("E38").Value = '  -2.21%  '

$ws.Range("E39").Value = '  -1.49%  '

$ws.Range("B40").Value = 'MultiversX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D40").Value = '64.66'
$ws.Range("E40").Value = '  +5.80%  '

$ws.Range("B41").Value = 'THORChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D41").Value = '5.66'
$ws.Range("E41").Value = '  -1.94%  '

$ws.Range("D42").Value = '0.200'
$ws.Range("E42").Value = '  -1.15%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '8.72'
$ws.Range("E43").Value = '  +0.66%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '104.71'
$ws.Range("E44").Value = '  -2.25%  '

$ws.Range("E45").Value = '  +1.16%  '

$ws.Range("E46").Value = '  +6.32%  '

$ws.Range("E47").Value = '  -0.20%  '

$ws.Range("D48").Value = '1.16'
$ws.Range("E48").Value = '  +0.43%  '

$ws.Range("D49").Value = '2.71'
$ws.Range("E49").Value = '  +0.76%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.431.73'
$ws.Range("E50").Value = '  -0.83%  '

$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").Value = '4.03'
$ws.Range("E51").Value = '  -1.03%  '
